# a69_f21_b UPP Pachuca - actualizacion de febrero (4to trimestre 2021)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet "Reporte de Formatos" (sheet1) - row 8 is the single data record.
# Period moves from Q3 2021 (01/07/2021-30/09/2021) to Q4 2021
# (01/10/2021-31/12/2021); validation/update dates move from 11/10/2021 to
# 10/01/2022; the hyperlink cell is cleared and the "Clasificacion ..." cell
# is cleared too; a note is added explaining the late / definitive data.
# ---------------------------------------------------------------------------

# Fecha de inicio / termino del periodo que se informa
$ws1.Cells.Item(8, 2).Value = 44470
$ws1.Cells.Item(8, 3).Value = 44561

# "Clasificacion del estado analitico ..." (D8) no longer filled in
$ws1.Cells.Item(8, 4).ClearContents()

# Drop the hyperlink + text that used to live in E8 ("Hipervinculo al Estado
# analitico ..."), then restore a plain (non-hyperlink) bordered look by
# reusing D8's existing format.
$ws1.Range("E8").Hyperlinks.Delete()
$ws1.Range("D8").Copy()
$ws1.Range("E8").PasteSpecial(-4122)
$ws1.Cells.Item(8, 5).ClearContents()
$excel.CutCopyMode = 0

# The workbook no longer uses the built-in "Hipervinculo" cell style once the
# hyperlink is gone - drop it (mirrors Excel dropping the unused style/font).
$wb.Styles.Item(1).Delete()

# Fecha de validacion / actualizacion
$ws1.Cells.Item(8, 7).Value = 44571
$ws1.Cells.Item(8, 8).Value = 44571

# Nota (I8): new justified, wrapped note text
$notaCell = $ws1.Cells.Item(8, 9)
$notaCell.Value = "Derivado del Cierre financiero - presupuestal que se trabaja de manera  coordinada entre planeación y administración para la entrega y preparación de información ante las diversas dependencias fiscalizadoras , en apego al artículo 15 de la ley de fiscalización superior y rendición de cuenta de estado de hidalgo,  correlativamente con la fracción V del artículo 28 de la misma ley se establecen como fecha de entrega  los siete días hábiles siguientes al cierre del trimestre.  Así también dentro del convenio especifico para la  asignación de recursos con carácter de apoyo solidario firmado entre la federación y el estado en su cláusula sexta fracción ""f"" donde obliga la entrega de los estados financieros dentro de los primeros diez días hábiles a la coordinación de universidades tecnológicas y politécnicas. Motivo por el cual estaremos entregando la información con cifras definitivas, el 28 de enero del año en curso."
$notaCell.HorizontalAlignment = -4130
$notaCell.WrapText = $true

# Row / column sizing for sheet1
$ws1.Rows.Item(3).RowHeight = 58.5
$ws1.Rows.Item(5).UseStandardHeight = $true
$ws1.Rows.Item(8).RowHeight = 203.25
$ws1.Columns.Item(5).ColumnWidth = 61.42578125
$ws1.Columns.Item(9).ColumnWidth = 73.140625

# ---------------------------------------------------------------------------
# Sheet "Tabla_393859" (sheet2) - the budget breakdown rows (chapters
# 1000-5000) are removed, keeping only the header rows.
# ---------------------------------------------------------------------------
$ws2.Range("A4:I8").EntireRow.Delete()
$ws2.Columns.Item(3).ColumnWidth = 38.5703125
$ws2.Columns.Item(6).ColumnWidth = 12.5703125
$ws2.Columns.Item(7).ColumnWidth = 12.85546875
$ws2.Columns.Item(8).ColumnWidth = 8.85546875

# Reset sheet2's lingering selection, then restore the real selection/active
# sheet (sheet1, cell H8) so tabSelected ends on the right sheet.
$ws2.Range("A1").Select()
$ws1.Range("H8").Select()
